$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.278.21'
$ws.Range('E2').Value = '  +3.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.816.14'
$ws.Range('E3').Value = '  +4.13%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '328.34'
$ws.Range('E5').Value = '  +2.09%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4353'
$ws.Range('E7').Value = '  +3.75%  '
$ws.Range('E8').Value = '  +2.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.02'
$ws.Range('E9').Value = '  -0.97%  '
$ws.Range('E10').Value = '  +3.74%  '
$ws.Range('E11').Value = '  +2.82%  '
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('E13').Value = '  +3.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.303'
$ws.Range('E14').Value = '  +3.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.538'
$ws.Range('E15').Value = '  +5.00%  '
$ws.Range('E16').Value = '  +4.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.27'
$ws.Range('E17').Value = '  +5.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001081'
$ws.Range('E18').Value = '  +1.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06539'
$ws.Range('E19').Value = '  +7.00%  '
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.51'
$ws.Range('E21').Value = '  +4.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.267'
$ws.Range('E22').Value = '  +3.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.300.78'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.64'
$ws.Range('E24').Value = '  +1.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.030'
$ws.Range('E25').Value = '  -13.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.95'
$ws.Range('E26').Value = '  +5.99%  '
$ws.Range('E27').Value = '  +2.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.033.71'
$ws.Range('E28').Value = '  +4.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.303'
$ws.Range('E29').Value = '  -2.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.90'
$ws.Range('E30').Value = '  +2.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.215'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.952'
$ws.Range('E32').Value = '  +5.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09189'
$ws.Range('E33').Value = '  +0.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.486'
$ws.Range('E34').Value = '  -4.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.93'
$ws.Range('E35').Value = '  +2.86%  '
$ws.Range('E36').Value = '  +2.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2180'
$ws.Range('E37').Value = '  +2.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.198'
$ws.Range('E38').Value = '  +2.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6574'
$ws.Range('E39').Value = '  +3.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06209'
$ws.Range('E40').Value = '  +2.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.193'
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.125'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.426'
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9997'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.91'
$ws.Range('E45').Value = '  +2.16%  '
$ws.Range('E46').Value = '  +4.91%  '
$ws.Range('E47').Value = '  +1.20%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.022'
$ws.Range('E48').Value = '  +4.24%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '125.65'
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.157'
$ws.Range('E50').Value = '  +4.08%  '
$ws.Range('E51').Value = '  +2.43%  '
